# Moved black and white to pin 8 and 9
# Re-positions the black/white pin-indicator ovals (and the connectors,
# battery picture and "Buttons" label that move together with them) on
# slide 2 of the Raspberry Pi Pico pinout diagram.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# EMU -> point conversion (1 pt = 12700 EMU), matching PowerPoint's
# Shape.Left / Shape.Top units. PowerPoint stores these as single-precision
# floats internally, and converting back to EMU on save truncates rather
# than rounds - so a plain division can land one EMU short. Nudge the
# point value up by tiny increments until the float32 round-trip reproduces
# the exact target EMU value.
$emuPerPt = 12700.0

function EmuToPt($targetEmu) {
    $base = $targetEmu / $emuPerPt
    for ($k = 0; $k -lt 2000; $k++) {
        $candidate = $base + ($k * 0.000001)
        $f32 = [float]$candidate
        $emu = [int64]($f32 * $emuPerPt)
        if ($emu -eq $targetEmu) {
            return $candidate
        }
    }
    return $base
}

# Shape 10: Oval 26 (black indicator oval)
$sh = $s.Shapes.Item(10)
$sh.Left = EmuToPt 1817452
$sh.Top  = EmuToPt 3650379

# Shape 11: Oval 27 (white indicator oval)
$sh = $s.Shapes.Item(11)
$sh.Left = EmuToPt 1817452
$sh.Top  = EmuToPt 3397015

# Shape 20: Straight Connector 41 (line into black oval's pin)
$sh = $s.Shapes.Item(20)
$sh.Left = EmuToPt 2023192
$sh.Top  = EmuToPt 3495599

# Shape 21: Straight Connector 43 (line into white oval's pin)
$sh = $s.Shapes.Item(21)
$sh.Left = EmuToPt 2023192
$sh.Top  = EmuToPt 3748963

# Shape 22: Picture 49 (battery icon)
$sh = $s.Shapes.Item(22)
$sh.Left = EmuToPt 1126238
$sh.Top  = EmuToPt 2373526

# Shape 23: Straight Connector 51 (battery lead line)
$sh = $s.Shapes.Item(23)
$sh.Left = EmuToPt 1817452
$sh.Top  = EmuToPt 2709585

# Shape 24: TextBox 52 ("Buttons" label)
$sh = $s.Shapes.Item(24)
$sh.Left = EmuToPt 1078756
$sh.Top  = EmuToPt 3975777
